# Updating SDTM annotations on draft Collection DSS
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# FTCAT row: annotate with the prepopulated category value
$ws.Range("AE4").Value = "FTCAT = SIX MINUTE WALK"

# FTORRESU rows: clarify that the unit is fixed to meters (m)
$ws.Range("AE7").Value  = "FTORRESU = m when FTTESTCD = SIXMW101"
$ws.Range("AE9").Value  = "FTORRESU = m when FTTESTCD = SIXMW102"
$ws.Range("AE11").Value = "FTORRESU = m when FTTESTCD = SIXMW103"
$ws.Range("AE13").Value = "FTORRESU = m when FTTESTCD = SIXMW104"
$ws.Range("AE15").Value = "FTORRESU = m when FTTESTCD = SIXMW105"
$ws.Range("AE17").Value = "FTORRESU = m when FTTESTCD = SIXMW106"

# Reflect the active cell selection at time of save
[void]$ws.Range("AE18").Select()
